$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2919
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2919
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2919
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3887
$ws.Range("H80").Value = 958
$ws.Range("I80").Value = 705.7143
$ws.Range("K80").Value = 2117.1429
$ws.Range("M80").Value = -1119.1429
$ws.Range("H83").Value = 958
$ws.Range("I83").Value = 705.7143
$ws.Range("K83").Value = 6351.428699999999
$ws.Range("M83").Value = -1359.428699999999
$ws.Range("H103").Value = 2364.1
$ws.Range("I103").Value = 2691.8572
$ws.Range("J103").Value = 1599.3334
$ws.Range("K103").Value = 8075.571599999999
$ws.Range("L103").Value = 4798.0002
$ws.Range("M103").Value = -7489.571599999999
$ws.Range("N103").Value = -5970.0002
$ws.Range("H111").Value = 2000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 6000
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -12134
$ws.Range("H112").Value = 2775.5557
$ws.Range("J112").Value = 2810
$ws.Range("L112").Value = 8430
$ws.Range("N112").Value = -10646
$ws.Range("H132").Value = 43315
$ws.Range("I132").Value = 43315
$ws.Range("K132").Value = 129945
$ws.Range("M132").Value = -127415
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 709.2
$ws.Range("I50").Value = 515.6667
$ws.Range("J50").Value = 999.5
$ws.Range("K50").Value = 515.6667
$ws.Range("L50").Value = 999.5
$ws.Range("M50").Value = 198.3333
$ws.Range("N50").Value = -2427.5
$ws.Range("H74").Value = 2203.1
$ws.Range("I74").Value = 2286.652
$ws.Range("K74").Value = 2286.652
$ws.Range("M74").Value = -1412.652
$ws.Range("H77").Value = 2203.1
$ws.Range("I77").Value = 2286.652
$ws.Range("K77").Value = 11433.26
$ws.Range("M77").Value = -7065.26
$ws.Range("H110").Value = 620.5833
$ws.Range("I110").Value = 625.3
$ws.Range("K110").Value = 625.3
$ws.Range("M110").Value = 1419.7
$ws.Range("H122").Value = 2826.077
$ws.Range("I122").Value = 2708.25
$ws.Range("K122").Value = 8124.75
$ws.Range("M122").Value = -5674.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5788.45
$ws.Range("I86").Value = 4143.364
$ws.Range("J86").Value = 7799.1113
$ws.Range("K86").Value = 4143.364
$ws.Range("L86").Value = 7799.1113
$ws.Range("M86").Value = -3020.364
$ws.Range("N86").Value = -10045.1113
$ws.Range("H89").Value = 5788.45
$ws.Range("I89").Value = 4143.364
$ws.Range("J89").Value = 7799.1113
$ws.Range("K89").Value = 20716.82
$ws.Range("L89").Value = 38995.5565
$ws.Range("M89").Value = -15100.82
$ws.Range("N89").Value = -50227.5565
$ws.Range("H105").Value = 2653.8823
$ws.Range("I105").Value = 2599.6667
$ws.Range("J105").Value = 2784
$ws.Range("K105").Value = 2599.6667
$ws.Range("L105").Value = 2784
$ws.Range("M105").Value = -852.6667000000002
$ws.Range("N105").Value = -6278
$ws.Range("H107").Value = 5941.4
$ws.Range("I107").Value = 5747.25
$ws.Range("J107").Value = 6070.8335
$ws.Range("K107").Value = 5747.25
$ws.Range("L107").Value = 6070.8335
$ws.Range("M107").Value = -3827.25
$ws.Range("N107").Value = -9910.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 24666.334
$ws.Range("J57").Value = 24666.334
$ws.Range("L57").Value = 24666.334
$ws.Range("N57").Value = -25786.334
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H134").Value = 1098.1666
$ws.Range("I134").Value = 1098.1666
$ws.Range("K134").Value = 3294.4998
$ws.Range("M134").Value = -759.4998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1964.9166
$ws.Range("I68").Value = 1903.4
$ws.Range("J68").Value = 2008.8572
$ws.Range("K68").Value = 5710.200000000001
$ws.Range("L68").Value = 6026.571599999999
$ws.Range("M68").Value = -4899.200000000001
$ws.Range("N68").Value = -7648.571599999999
$ws.Range("H71").Value = 1964.9166
$ws.Range("I71").Value = 1903.4
$ws.Range("J71").Value = 2008.8572
$ws.Range("K71").Value = 17130.6
$ws.Range("L71").Value = 18079.7148
$ws.Range("M71").Value = -13074.6
$ws.Range("N71").Value = -26191.7148
$ws.Range("H103").Value = 676.5
$ws.Range("I103").Value = 246.4
$ws.Range("K103").Value = 739.2
$ws.Range("M103").Value = 139.8
$ws.Range("H107").Value = 332
$ws.Range("I107").Value = 330.66666
$ws.Range("J107").Value = 336
$ws.Range("K107").Value = 991.9999799999999
$ws.Range("L107").Value = 1008
$ws.Range("M107").Value = 928.0000200000001
$ws.Range("N107").Value = -4848
$ws.Range("H113").Value = 492.84616
$ws.Range("I113").Value = 540.8
$ws.Range("J113").Value = 333
$ws.Range("K113").Value = 1622.4
$ws.Range("L113").Value = 999
$ws.Range("M113").Value = 547.6000000000001
$ws.Range("N113").Value = -5339
$ws.Range("H125").Value = 3950
$ws.Range("I125").Value = 3950
$ws.Range("K125").Value = 11850
$ws.Range("M125").Value = -6930
$ws.Range("H132").Value = 1841.2941
$ws.Range("I132").Value = 1789.7
$ws.Range("J132").Value = 1915
$ws.Range("K132").Value = 16107.3
$ws.Range("L132").Value = 17235
$ws.Range("M132").Value = -13577.3
$ws.Range("N132").Value = -22295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 1540.0769
$ws.Range("I102").Value = 1662.2
$ws.Range("K102").Value = 1662.2
$ws.Range("M102").Value = -40.20000000000005
$ws.Range("H123").Value = 975000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2788.5833
$ws.Range("I55").Value = 1280.4286
$ws.Range("J55").Value = 4900
$ws.Range("K55").Value = 1280.4286
$ws.Range("L55").Value = 4900
$ws.Range("M55").Value = -1107.4286
$ws.Range("N55").Value = -5246
$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251
$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256
$ws.Range("H82").Value = 7142.857
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 7142.857
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H93").Value = 1200
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H122").Value = 3997
$ws.Range("I122").Value = 3997
$ws.Range("K122").Value = 11991
$ws.Range("M122").Value = -9541

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 14997.5
$ws.Range("J50").Value = 14997.5
$ws.Range("L50").Value = 14997.5
$ws.Range("N50").Value = -16259.5
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820
